$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 on the "Rules" sheet held the shared string "R40"; it is replaced
# with the literal text "1" (stored as text, not a number -- force the
# number format to Text first so the COM layer keeps it a string value).
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
